# Appointment_List: add a new appointment row (row 5) to the sheet.
# Columns: A=Appointment ID, B=Patient ID, C=Doctor ID, D=Status,
#          E=Appointment Date, F=Appointment Time, G=Outcome Record
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "test"
$ws.Range("B5").Value = "P1001"
$ws.Range("C5").Value = "D001"
$ws.Range("D5").Value = "SCHEDULED"
$ws.Range("E5").Value = 45595
$ws.Range("F5").Value = "9:00 am"

# Outcome Record starts out blank for a freshly-created appointment.
# Writing an apostrophe-prefixed empty string forces Excel to store an
# (empty) text entry in the cell rather than clearing it outright, then
# resetting the style keeps the cell formatting the same as its neighbours.
$ws.Range("G5").Value = "'"
$ws.Range("G5").Style = "Normal"
